$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the longer matchup descriptions
$ws.Columns.Item(1).ColumnWidth = 49.63

# Add new rows of match results
$ws.Range("A5").Value = "Random++ vs Random+"
$ws.Range("B5").Value = 51
$ws.Range("C5").Value = 8

$ws.Range("A6").Value = "Random++ vs Random++"
$ws.Range("B6").Value = 37
$ws.Range("C6").Value = 51

$ws.Range("A7").Value = "Alpha 10 samplesize vs Random++"
$ws.Range("B7").Value = 51
$ws.Range("C7").Value = 3

$ws.Range("A8").Value = "Alpha 50 vs Random++ "
$ws.Range("B8").Value = 51
$ws.Range("C8").Value = 0

$ws.Range("A9").Value = "Alpha 10 vs Alpha 50"
$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 51

# Move the selection to match the saved view state
$ws.Range("D9").Select()

$wb.Save()
